# Applies the commit:
#  1) Re-points the three "Table_0"-styled tables (slides 14-16) at the
#     built-in table style {3777AE3A-40FC-4C1E-88CC-073C0E44CC49}.
#  2) Swaps the slide-master theme's colour scheme from the custom
#     "Integral / Red Violet" palette to the stock "Office" palette
#     (the net effect, for the part of the theme swap that is reachable
#     through the exposed PowerPoint object model, of the theme1.xml /
#     theme2.xml content swap recorded in the commit).

$p = $ppt.ActivePresentation

# --- 1) Retarget table styles -------------------------------------------
$oldStyleId = "{5179044E-BC99-43DE-8D43-5142F525A787}"
$newStyleId = "{3777AE3A-40FC-4C1E-88CC-073C0E44CC49}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2) Swap the theme colour scheme on the slide master ----------------
# PowerPoint's RGB property packs colours as 0x00BBGGRR (Windows COLORREF),
# i.e. byte-reversed relative to the "RRGGBB" hex used in the OOXML
# <a:srgbClr val="RRGGBB"/> markup.
function BGR([int]$r, [int]$g, [int]$b) {
    return $b * 65536 + $g * 256 + $r
}

# Target palette: the stock "Office" theme colours (was "Red Violet").
$officeColors = @(
    (BGR 0x00 0x00 0x00),  # 1  dk1
    (BGR 0xFF 0xFF 0xFF),  # 2  lt1
    (BGR 0x44 0x54 0x6A),  # 3  dk2
    (BGR 0xE7 0xE6 0xE6),  # 4  lt2
    (BGR 0x5B 0x9B 0xD5),  # 5  accent1
    (BGR 0xED 0x7D 0x31),  # 6  accent2
    (BGR 0xA5 0xA5 0xA5),  # 7  accent3
    (BGR 0xFF 0xC0 0x00),  # 8  accent4
    (BGR 0x44 0x72 0xC4),  # 9  accent5
    (BGR 0x70 0xAD 0x47),  # 10 accent6
    (BGR 0x05 0x63 0xC1),  # 11 hlink
    (BGR 0x95 0x4F 0x72)   # 12 folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($k = 1; $k -le $themeColors.Count; $k++) {
    $themeColors.Colors($k).RGB = $officeColors[$k - 1]
}
